$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("DATA")
$ws2 = $wb.Worksheets.Item("OUTPUT")
$ws3 = $wb.Worksheets.Item("CALC")

# Fix room overloading (not a sort event): swap the date-group stamp
# ("101012020" / "101012021") between the first row and the block of
# rows that had been mistakenly tagged with the other group's value.
$ws1.Range("G3").Value = 101012021

$ws1.Range("G84").Value = 101012020
$ws1.Range("G85").Value = 101012020
$ws1.Range("G86").Value = 101012020
$ws1.Range("G87").Value = 101012020
$ws1.Range("G88").Value = 101012020
$ws1.Range("G89").Value = 101012020
$ws1.Range("G91").Value = 101012020
$ws1.Range("G92").Value = 101012020
$ws1.Range("G105").Value = 101012020
$ws1.Range("G109").Value = 101012020
$ws1.Range("G111").Value = 101012020
$ws1.Range("G113").Value = 101012020
$ws1.Range("G114").Value = 101012020
$ws1.Range("G122").Value = 101012020
$ws1.Range("G123").Value = 101012020
$ws1.Range("G124").Value = 101012020
$ws1.Range("G125").Value = 101012020

# Restore the UI state the author left the workbook in when saving:
# DATA's cursor parked on E25 (not the active tab any more), OUTPUT
# as the active/selected tab with its cursor on J12 (CALC, previously
# active, naturally loses tabSelected once OUTPUT is selected last).
$ws1.Range("E25").Select()
$ws2.Range("J12").Select()
